# Basic excel file processing and ETL
# Rename the header row labels to friendlier, spaced-out captions and
# move the active selection to reflect where the user was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row relabel (row 1, columns A-D)
$ws.Range("A1").Value = "Effective Date"
$ws.Range("B1").Value = "Client Account"
$ws.Range("C1").Value = "External Reference"
$ws.Range("D1").Value = "Company Name"

# Move/save the active cell selection on the sheet
$ws.Range("E12").Select()
